$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 27784514
$ws.Range("I62").Value = 1377.0714
$ws.Range("J62").Value = 125025496
$ws.Range("K62").Value = 1377.0714
$ws.Range("L62").Value = 125025496
$ws.Range("M62").Value = -753.0714
$ws.Range("N62").Value = -125026744

$ws.Range("H63").Value = 314128.34
$ws.Range("J63").Value = 314128.34
$ws.Range("L63").Value = 314128.34
$ws.Range("N63").Value = -315376.34

$ws.Range("H65").Value = 27784514
$ws.Range("I65").Value = 1377.0714
$ws.Range("J65").Value = 125025496
$ws.Range("K65").Value = 6885.357
$ws.Range("L65").Value = 625127480
$ws.Range("M65").Value = -3765.357
$ws.Range("N65").Value = -625133720

$ws.Range("H66").Value = 314128.34
$ws.Range("J66").Value = 314128.34
$ws.Range("L66").Value = 942385.02
$ws.Range("N66").Value = -948625.02

$ws.Range("H80").Value = 1200
$ws.Range("J80").Value = 1200
$ws.Range("L80").Value = 3600
$ws.Range("N80").Value = -5596

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H83").Value = 1200
$ws.Range("J83").Value = 1200
$ws.Range("L83").Value = 10800
$ws.Range("N83").Value = -20784

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H86").Value = 11663.75
$ws.Range("I86").Value = 19067.666
$ws.Range("J86").Value = 4259.8335
$ws.Range("K86").Value = 19067.666
$ws.Range("L86").Value = 4259.8335
$ws.Range("M86").Value = -17944.666
$ws.Range("N86").Value = -6505.8335

$ws.Range("H89").Value = 11663.75
$ws.Range("I89").Value = 19067.666
$ws.Range("J89").Value = 4259.8335
$ws.Range("K89").Value = 95338.33
$ws.Range("L89").Value = 21299.1675
$ws.Range("M89").Value = -89722.33
$ws.Range("N89").Value = -32531.1675

$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

$ws.Range("H105").Value = 32042.75
$ws.Range("J105").Value = 32042.75
$ws.Range("L105").Value = 32042.75
$ws.Range("N105").Value = -39030.75

$ws.Range("H107").Value = 83929.336
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 1000
$ws.Range("N107").Value = -4840

$ws.Range("H110").Value = 42500
$ws.Range("J110").Value = 42500
$ws.Range("L110").Value = 42500
$ws.Range("N110").Value = -50680

$ws.Range("H112").Value = 2317.5862
$ws.Range("J112").Value = 2826.5217
$ws.Range("L112").Value = 8479.5651
$ws.Range("N112").Value = -10695.5651

$ws.Range("H113").Value = 26319468
$ws.Range("I113").Value = 62501820
$ws.Range("J113").Value = 5030.727
$ws.Range("K113").Value = 62501820
$ws.Range("L113").Value = 5030.727
$ws.Range("M113").Value = -62498566
$ws.Range("N113").Value = -11538.727

$ws.Range("H114").Value = 36666.668
$ws.Range("J114").Value = 36666.668
$ws.Range("L114").Value = 36666.668
$ws.Range("N114").Value = -45344.668

$ws.Range("H116").Value = 4504.552
$ws.Range("I116").Value = 2387.5
$ws.Range("J116").Value = 7110.154
$ws.Range("K116").Value = 2387.5
$ws.Range("L116").Value = 7110.154
$ws.Range("M116").Value = 1054.5
$ws.Range("N116").Value = -13994.154

$ws.Range("H117").Value = 35000
$ws.Range("J117").Value = 35000
$ws.Range("L117").Value = 35000
$ws.Range("N117").Value = -44178

$ws.Range("H132").Value = 3557.0322
$ws.Range("I132").Value = 2517.524
$ws.Range("J132").Value = 5740
$ws.Range("K132").Value = 7552.572
$ws.Range("L132").Value = 17220
$ws.Range("M132").Value = -5022.572
$ws.Range("N132").Value = -22280

$ws.Range("H135").Value = 34796.3
$ws.Range("I135").Value = 51367.25
$ws.Range("J135").Value = 1654.4
$ws.Range("K135").Value = 462305.25
$ws.Range("L135").Value = 14889.6
$ws.Range("M135").Value = -459770.25
$ws.Range("N135").Value = -19959.6

$ws.Range("H137").Value = 1626.6333
$ws.Range("I137").Value = 1288.909
$ws.Range("J137").Value = 2555.375
$ws.Range("K137").Value = 3866.727
$ws.Range("L137").Value = 7666.125
$ws.Range("M137").Value = -1316.727
$ws.Range("N137").Value = -12766.125

$ws.Range("H141").Value = 2454.682
$ws.Range("I141").Value = 1780.875
$ws.Range("J141").Value = 4251.5
$ws.Range("K141").Value = 5342.625
$ws.Range("L141").Value = 12754.5
$ws.Range("M141").Value = -162.625
$ws.Range("N141").Value = -23114.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2935.1724
$ws.Range("I2").Value = 1072
$ws.Range("J2").Value = 7075.5557
$ws.Range("K2").Value = 1072
$ws.Range("L2").Value = 7075.5557
$ws.Range("M2").Value = -959
$ws.Range("N2").Value = -7301.5557

$ws.Range("H61").Value = 2348.6
$ws.Range("I61").Value = 2360.75
$ws.Range("K61").Value = 2360.75
$ws.Range("M61").Value = -2148.75

$ws.Range("H116").Value = 2935.1724
$ws.Range("I116").Value = 1072
$ws.Range("J116").Value = 7075.5557
$ws.Range("K116").Value = 1072
$ws.Range("L116").Value = 7075.5557
$ws.Range("M116").Value = 1222
$ws.Range("N116").Value = -11663.5557

$ws.Range("H136").Value = 2348.6
$ws.Range("I136").Value = 2360.75
$ws.Range("K136").Value = 7082.25
$ws.Range("M136").Value = -4532.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2935.1724
$ws.Range("I3").Value = 1072
$ws.Range("J3").Value = 7075.5557
$ws.Range("K3").Value = 1072
$ws.Range("L3").Value = 7075.5557
$ws.Range("M3").Value = -958
$ws.Range("N3").Value = -7303.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2615.7
$ws.Range("I132").Value = 2281.375
$ws.Range("J132").Value = 3953
$ws.Range("K132").Value = 6844.125
$ws.Range("L132").Value = 11859
$ws.Range("M132").Value = -4314.125
$ws.Range("N132").Value = -16919

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 499.82053
$ws.Range("I5").Value = 317.44116
$ws.Range("J5").Value = 1740
$ws.Range("K5").Value = 952.32348
$ws.Range("L5").Value = 5220
$ws.Range("M5").Value = -840.32348
$ws.Range("N5").Value = -5444

$ws.Range("H68").Value = 1077.0769
$ws.Range("I68").Value = 917
$ws.Range("J68").Value = 1214.2858
$ws.Range("K68").Value = 2751
$ws.Range("L68").Value = 3642.8574
$ws.Range("M68").Value = -1940
$ws.Range("N68").Value = -5264.857400000001

$ws.Range("H71").Value = 1077.0769
$ws.Range("I71").Value = 917
$ws.Range("J71").Value = 1214.2858
$ws.Range("K71").Value = 8253
$ws.Range("L71").Value = 10928.5722
$ws.Range("M71").Value = -4197
$ws.Range("N71").Value = -19040.5722

$ws.Range("H107").Value = 1572
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 2153.3333
$ws.Range("K107").Value = 2100
$ws.Range("L107").Value = 6459.999899999999
$ws.Range("M107").Value = -180
$ws.Range("N107").Value = -10299.9999

$ws.Range("H113").Value = 543.4318
$ws.Range("I113").Value = 536.8823
$ws.Range("J113").Value = 547.55554
$ws.Range("K113").Value = 1610.6469
$ws.Range("L113").Value = 1642.66662
$ws.Range("M113").Value = 559.3531
$ws.Range("N113").Value = -5982.66662

$ws.Range("H122").Value = 45702.914
$ws.Range("I122").Value = 255.28572
$ws.Range("J122").Value = 49359.62
$ws.Range("K122").Value = 2297.57148
$ws.Range("L122").Value = 444236.58
$ws.Range("M122").Value = 152.4285199999999
$ws.Range("N122").Value = -449136.58

$ws.Range("H132").Value = 1097517.9
$ws.Range("I132").Value = 1197019.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 10773175.5
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -10770645.5
$ws.Range("N132").Value = -32060

$ws.Range("H133").Value = 7998.5713
$ws.Range("I133").Value = 7632.5
$ws.Range("J133").Value = 8486.666999999999
$ws.Range("K133").Value = 22897.5
$ws.Range("L133").Value = 25460.001
$ws.Range("M133").Value = -17837.5
$ws.Range("N133").Value = -35580.001

$ws.Range("H134").Value = 2532.7273
$ws.Range("I134").Value = 1710
$ws.Range("J134").Value = 3972.5
$ws.Range("K134").Value = 5130
$ws.Range("L134").Value = 11917.5
$ws.Range("M134").Value = -60
$ws.Range("N134").Value = -22057.5

$ws.Range("H135").Value = 499.82053
$ws.Range("I135").Value = 317.44116
$ws.Range("J135").Value = 1740
$ws.Range("K135").Value = 2856.97044
$ws.Range("L135").Value = 15660
$ws.Range("M135").Value = -321.9704400000001
$ws.Range("N135").Value = -20730

$ws.Range("H138").Value = 47621210
$ws.Range("I138").Value = 166667800
$ws.Range("J138").Value = 2568.7334
$ws.Range("K138").Value = 500003400
$ws.Range("L138").Value = 7706.2002
$ws.Range("M138").Value = -499998260
$ws.Range("N138").Value = -17986.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3748.75
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4497.5
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 13492.5
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -18392.5

$ws.Range("H132").Value = 3321.1428
$ws.Range("I132").Value = 2099.7778
$ws.Range("K132").Value = 6299.3334
$ws.Range("M132").Value = -3769.3334

$ws.Range("H136").Value = 1794.85
$ws.Range("I136").Value = 1626.1578
$ws.Range("K136").Value = 4878.4734
$ws.Range("M136").Value = -2328.4734

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1467.2941
$ws.Range("I136").Value = 1297.125
$ws.Range("J136").Value = 4190
$ws.Range("K136").Value = 3891.375
$ws.Range("L136").Value = 12570
$ws.Range("M136").Value = -1341.375
$ws.Range("N136").Value = -17670

$ws.Range("H139").Value = 39800
$ws.Range("J139").Value = 39800
$ws.Range("L139").Value = 39800
$ws.Range("N139").Value = -50080
